$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new day column "06-nov" right after
# column DJ ("05-nov"), pushing the "01-oct." ... "31-oct." block one
# column to the right (DK:EO -> DL:EP). The freshly inserted column gets
# the same "-" placeholder used for other not-yet-available days.
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("DK1").EntireColumn.Insert()
$wsPrix.Range("DK1").Value = "06-nov"

for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 115).Value = "-"
}

# --- Sheet "Gaz": append the next day's price row.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Cells.Item(143, 1).NumberFormat = "@"
$wsGaz.Cells.Item(143, 1).Value = "2025-11-04"
$wsGaz.Cells.Item(143, 1).Style = "Normal"
$wsGaz.Cells.Item(143, 2).Value = 31.17

# --- Sheet "CO2": append the next day's price row.
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Cells.Item(143, 1).NumberFormat = "@"
$wsCo2.Cells.Item(143, 1).Value = "2025-11-04"
$wsCo2.Cells.Item(143, 1).Style = "Normal"
$wsCo2.Cells.Item(143, 2).Value = 81.9
